$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 22

$ws.Cells.Item($row, 1).Value = 7
$ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value = "Ñuble"
$ws.Cells.Item($row, 4).Value = 45166
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat
$ws.Cells.Item($row, 5).Value = 16
$ws.Cells.Item($row, 6).Value = 100112039
$ws.Cells.Item($row, 7).Value = "Ciboulette"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 120
$ws.Cells.Item($row, 11).Value = 2500
$ws.Cells.Item($row, 12).Value = 2500
$ws.Cells.Item($row, 13).Value = 2500
$ws.Cells.Item($row, 14).Value = "$/docena de atados"
$ws.Cells.Item($row, 15).Value = "Región Metropolitana"
$ws.Cells.Item($row, 16).Value = 833
$ws.Cells.Item($row, 17).Value = 3
$ws.Cells.Item($row, 18).Value = "Hortaliza"
